$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A5"))
$ws.Sort.SetRange($ws.Range("A2:B5"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

$ws.Range("A2:B5").Select()
